$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-12 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-05-13 Monday", 2)

$d.Content.Find.Execute("576×5=2880", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "135×2=270", 2)

$d.Content.Find.Execute("443×3=1329", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "579×2=1158", 2)

$d.Content.Find.Execute("870×6=5220", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "902×3=2706", 2)

$d.Content.Find.Execute("492×3=1476", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "726×2=1452", 2)

$d.Content.Find.Execute("909×2=1818", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "455×7=3185", 2)

$d.Content.Find.Execute("680×4=2720", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "419×9=3771", 2)

$d.Content.Find.Execute("149×9=1341", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "193×2=386", 2)

$d.Content.Find.Execute("181×8=1448", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "132×2=264", 2)

$d.Content.Find.Execute("628×7=4396", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "232×6=1392", 2)

$d.Content.Find.Execute("228×4=912", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "717×5=3585", 2)

$d.Content.Find.Execute("691×6=4146", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "826×2=1652", 2)

$d.Content.Find.Execute("517×2=1034", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "682×2=1364", 2)

$d.Content.Find.Execute("375×2=750", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "185×8=1480", 2)

$d.Content.Find.Execute("491×6=2946", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "623×4=2492", 2)

$d.Content.Find.Execute("644×8=5152", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "976×4=3904", 2)

$d.Content.Find.Execute("237×4=948", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "714×6=4284", 2)

$d.Content.Find.Execute("895×2=1790", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "486×6=2916", 2)

$d.Content.Find.Execute("291×9=2619", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "446×3=1338", 2)

$d.Content.Find.Execute("918×2=1836", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "838×5=4190", 2)

$d.Content.Find.Execute("448×5=2240", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "668×9=6012", 2)

$d.Content.Find.Execute("684×6=4104", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "958×7=6706", 2)

$d.Content.Find.Execute("357×3=1071", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "791×5=3955", 2)

$d.Content.Find.Execute("934×9=8406", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "898×6=5388", 2)

$d.Content.Find.Execute("224×5=1120", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "830×7=5810", 2)

$d.Content.Find.Execute("473×3=1419", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "540×7=3780", 2)

Write-Output "Done"
